# Apply cryptos list update (prices & 1h volume percentages)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.780.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "'1.949.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'248.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.88%  "
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").Value = "'0.06829"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "'112.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.14%  "
$ws.Range("D11").Value = "'19.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.34%  "
$ws.Range("D12").Value = "'1.942.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "'5.571"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.82%  "
$ws.Range("D14").Value = "'0.07659"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'0.6936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "'298.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.58%  "
$ws.Range("D17").Value = "'30.812.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'13.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.22%  "
$ws.Range("D19").Value = "'5.714"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.31%  "
$ws.Range("D20").Value = "'0.000007707"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "'2.206.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'0.9994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'6.597"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").Value = "'9.751"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.64%  "
$ws.Range("D26").Value = "'168.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").Value = "'20.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").Value = "'2.183"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("D29").Value = "'0.1092"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Value = "'4.635"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.17%  "
$ws.Range("D32").Value = "'4.446"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.53%  "
$ws.Range("D33").Value = "'0.05080"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'0.7796"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.23%  "
$ws.Range("D35").Value = "'1.167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("D36").Value = "'0.02076"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("D37").Value = "'2.732"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'2.707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").Value = "'2.046"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").Value = "'111.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "'0.4483"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").Value = "'0.8778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "'71.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("D45").Value = "'0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'7.429"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").Value = "'9.520"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "'49.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").Value = "'0.1262"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Value = "'35.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("D51").Value = "'0.2564"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.16%  "
